$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the "Mvc" paragraph (new content is inserted right before it).
# ------------------------------------------------------------------
$mvcIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Mvc") {
        $mvcIndex = $i
        break
    }
}
if ($mvcIndex -eq -1) {
    throw "Could not locate the 'Mvc' paragraph"
}

$mvcPara = $d.Paragraphs.Item($mvcIndex)

# Create an empty placeholder paragraph immediately before "Mvc", then
# fill it with the eight new paragraphs via a single InsertXML call
# (keeps each new paragraph's own pPr/rPr/proofErr markup intact).
$mvcPara.Range.InsertParagraphBefore()

$placeholder = $d.Paragraphs.Item($mvcIndex)

$newXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Ajax for number of users in admin page</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SEO</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Blog with comments inserting. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Assignments code to make one more </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>feature</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Search </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>on</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the website for elements</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Subscription emails giving emails to the </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>address.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$placeholder.Range.InsertXML($newXml)

# ------------------------------------------------------------------
# 2. Remove the old "Ajax for number of classes" / "Search in the
#    website for elements" / "Subscription emails giving emails to
#    the adress" / "Blog with comments inserting." paragraphs,
#    collapsing them down into a single empty paragraph.
# ------------------------------------------------------------------
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Ajax for number of classes") {
        $startIndex = $i
        break
    }
}
if ($startIndex -eq -1) {
    throw "Could not locate the 'Ajax for number of classes' paragraph"
}

$endIndex = -1
for ($i = $startIndex; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Blog with comments inserting.") {
        $endIndex = $i
        break
    }
}
if ($endIndex -eq -1) {
    throw "Could not locate the 'Blog with comments inserting.' paragraph"
}

$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($endIndex)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

Write-Output "done"
